# Pull in the 3rd (most recent) voucher/tender item.
# A new row is inserted at row 2 (pushing the existing rows down) and
# populated with the latest quotation details from the Ministry of
# Education (MOE000ETQ21000131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 - existing rows 2 & 3
# shift down to become rows 3 & 4.
$ws.Rows.Item(2).Insert()

# Populate the new row with the most recent tender.
$ws.Cells.Item(2, 1).Value2 = "Quotation - MOE000ETQ21000131"
$ws.Cells.Item(2, 2).Value2 = "INVITATION TO QUOTE FOR BOOK VOUCHERS"
$ws.Cells.Item(2, 3).Value2 = "19 May 2021`n01:00PM"
$ws.Cells.Item(2, 3).WrapText = $true
$ws.Cells.Item(2, 4).Value2 = "Ministry of Education"
$ws.Cells.Item(2, 5).Value2 = "Administration & Training ⇒ Gifts & Souvenirs"

# Match the row height used by the other data rows.
$ws.Rows.Item(2).RowHeight = 44.25

# Refresh column widths to fit the (now wider) content and move the
# selection onto the newly added row, mirroring the manual edit.
$ws.Columns.AutoFit()
$ws.Cells.Item(2, 1).Select()
